# Master file and variables list update — adding new wind/solar CAPEX/OPEX/
# load-factor input variables (WWIC/WOIC/WSIC, WWOM/WOOM/WSOM, WWLF/WOLF/WSLF)
# to the FTT-H2 variable list and the Time_Horizons lookup sheet.

$wb = $excel.ActiveWorkbook

$wsH2 = $wb.Worksheets.Item("FTT-H2")
$wsTH = $wb.Worksheets.Item("Time_Horizons")

# ---------------------------------------------------------------------------
# FTT-H2 sheet: 9 new variable rows (36-44), following the existing layout:
# A=Code, B=Read in?, C=numeric id, D=Description, E=RowDim, F=ColDim,
# G=3DDim, H=Conversion?, I=Scenario
# ---------------------------------------------------------------------------
$newVars = @(
    @{ Row = 36; Code = "WWIC"; Id = 36390000; Desc = "Onshore wind CAPEX" },
    @{ Row = 37; Code = "WOIC"; Id = 36400000; Desc = "Offshore wind CAPEX" },
    @{ Row = 38; Code = "WSIC"; Id = 36410000; Desc = "Solar CAPEX" },
    @{ Row = 39; Code = "WWOM"; Id = 36420000; Desc = "Onshore wind OPEX" },
    @{ Row = 40; Code = "WOOM"; Id = 36430000; Desc = "Offshore wind OPEX" },
    @{ Row = 41; Code = "WSOM"; Id = 36440000; Desc = "Solar OPEX" },
    @{ Row = 42; Code = "WWLF"; Id = 36450000; Desc = "Onshore wind load factor" },
    @{ Row = 43; Code = "WOLF"; Id = 36460000; Desc = "Offshore wind load factor" },
    @{ Row = 44; Code = "WSLF"; Id = 36470000; Desc = "Solar load factor" }
)

# Column A (variable codes) is filled in top-to-bottom first ...
foreach ($row in $newVars) {
    $r = $row.Row
    $wsH2.Cells.Item($r, 1).Value = $row.Code
    $wsH2.Cells.Item($r, 1).Font.Color = 0
}

# ... then the remaining columns, with column D (descriptions) entered in the
# same order the source workbook used (load-factor rows bottom-up, then the
# CAPEX/OPEX rows top-down) so new shared-string indices line up.
$descOrder = @(44, 43, 42, 36, 37, 38, 39, 40, 41)
foreach ($r in $descOrder) {
    $desc = ($newVars | Where-Object { $_.Row -eq $r }).Desc
    $wsH2.Cells.Item($r, 4).Value = $desc
    $wsH2.Cells.Item($r, 4).Font.Color = 0
}

foreach ($row in $newVars) {
    $r = $row.Row
    $wsH2.Cells.Item($r, 2).Value = 1
    $wsH2.Cells.Item($r, 3).Value = $row.Id
    $wsH2.Cells.Item($r, 5).Value = "RSHORTTI"
    $wsH2.Cells.Item($r, 6).Value = "TIME"
    $wsH2.Cells.Item($r, 7).Value = 0
    $wsH2.Cells.Item($r, 8).Value = 0
    $wsH2.Cells.Item($r, 9).Value = "S0"
}

# Three trailing blank (but styled) rows below the new data
$wsH2.Cells.Item(45, 1).Font.Color = 0
$wsH2.Cells.Item(46, 1).Font.Color = 0
$wsH2.Cells.Item(47, 1).Font.Color = 0

# ---------------------------------------------------------------------------
# Time_Horizons sheet: map each new variable code to its read-in horizon tag
# ---------------------------------------------------------------------------
$newCodes = @("WWIC", "WOIC", "WSIC", "WWOM", "WOOM", "WSOM", "WWLF", "WOLF", "WSLF")
$r = 101
foreach ($code in $newCodes) {
    $wsTH.Cells.Item($r, 1).Value = $code
    $wsTH.Cells.Item($r, 2).Value = "tl_2010"
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# View state: FTT-H2 becomes the active / selected sheet, with the last
# edited cells left selected on each touched sheet.
# ---------------------------------------------------------------------------
$wsTH.Range("B109").Select()
$wsH2.Activate()
$wsH2.Range("C44").Select()
